# Add team record (Wins/Losses/Ties) columns to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new columns AD="Wins", AE="Losses", AF="Ties" ---
# Copy the formatting of the existing header cell (AC1, style index 1: bold,
# centered, bordered) onto the three new header cells so they match the
# look of the rest of the header row.
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (2-50): every player/row gets the team's season record ---
$lastRow = 50
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Range("AD$r").Value = 83
    $ws.Range("AE$r").Value = 79
    $ws.Range("AF$r").Value = 0
}
